$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new "Test_column" before the existing "Percent_Employed"
# column (V). This pushes the old Percent_Employed data from column V to
# column W, leaving the new column V's data cells (V2:V9) blank.
$ws.Columns("V").Insert()
$ws.Range("V1").Value = "Test_column"

# Recompute the (moved) Percent_Employed column in W as Employed / Total,
# matching the "updated the census api work" refresh of that figure.
for ($row = 2; $row -le 9; $row++) {
    $employed = $ws.Range("J$row").Value2
    $total = $ws.Range("E$row").Value2
    $ws.Range("W$row").Value = $employed / $total
}

# Fill in the (previously blank) Major_category column with "Arts" for
# every data row.
$ws.Range("D2:D9").Value = "Arts"
